
# "Adding corrections to labs 1-3"
# Fix the three transposed/mis-typed values in the hemorrhage data table on
# Sheet1 and restore the cursor/selection to where the user left it
# (scrolled down to row 4, with F4:F5 selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Data corrections -------------------------------------------------
# Row 16 / Column C ("72-hour" column) corrected from 52 -> 84
$ws.Range("C16").Value = 84

# Row 18 / Column C corrected from 84 -> 52 (values had been swapped)
$ws.Range("C18").Value = 52

# Row 25 / Column C corrected from 0.13 -> 0.01
$ws.Range("C25").Value = 0.01

# --- View / selection state --------------------------------------------
# Leave the selection on F4:F5, scrolled so row 4 is at the top, matching
# where the editor was working when the corrections were made.
[void]$ws.Range("F4:F5").Select()
